## financeiraJavaFX scrum sheet update
## - A6 "GUI" -> "validaçoes"
## - New cell E6 "validações data e dinheiro (Como Classe de Eventos)"
##   (reuses the same fgColor/bgColor fill already used by A6/B9 - fillId 6)
## - New cell D7 "validações data e dinheiro" (same fill reuse)
## - New row 11: A11 "GUI" (moved out of A6), with a new light-yellow fill,
##   centered; B11 left blank
## - Column E widened
## - Selection moved to D8

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A6: "GUI" -> "validaçoes" (keep the existing fill/format as-is)
$ws.Range("A6").Value = "validaçoes"

# E6: new cell - copy A6's format (same centered fill) then set its text
$ws.Range("A6").Copy()
$ws.Range("E6").PasteSpecial(-4122)
$ws.Range("E6").Value = "validações data e dinheiro (Como Classe de Eventos)"

# D7: new cell - same format reuse
$ws.Range("A6").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D7").Value = "validações data e dinheiro"

$excel.CutCopyMode = 0

# New row 11: A11 "GUI" with a new light-yellow fill, centered; B11 blank
$ws.Range("A11").Value = "GUI"
$ws.Range("A11").HorizontalAlignment = -4108
$ws.Range("A11").Interior.Color = 13434879

# Widen column E (closest achievable value to the authored 44.78 "chars" width;
# the stored OOXML width is quantized to whole pixels by the column-width model)
$ws.Columns("E").ColumnWidth = 44.0

# Move the active selection to D8 (matches the author's last cursor position)
$ws.Range("D8").Select() | Out-Null
